$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A1").Value = 1.681432485580444
$ws.Range("B1").Value = 3.312968254089355
$ws.Range("C1").Value = 5.945855140686035
$ws.Range("D1").Value = 1.82273006439209
$ws.Range("E1").Value = 0.8992128968238831
